# Update crypto price list data (prices, 1h volume %, and two swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.430.97'
$ws.Cells.Item(2, 5).Value = '  +0.18%  '

$ws.Cells.Item(3, 4).Value = '1.848.47'
$ws.Cells.Item(3, 5).Value = '  +0.32%  '

$ws.Cells.Item(4, 5).Value = '  +0.16%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '240.86'
$ws.Cells.Item(5, 5).Value = '  +0.81%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.6277'
$ws.Cells.Item(6, 5).Value = '  -0.44%  '

$ws.Cells.Item(7, 5).Value = '  +0.09%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07678'
$ws.Cells.Item(8, 5).Value = '  +1.96%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.2919'
$ws.Cells.Item(9, 5).Value = '  -0.31%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '24.83'
$ws.Cells.Item(10, 5).Value = '  +1.60%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07746'

$ws.Cells.Item(12, 4).Value = '1.846.09'
$ws.Cells.Item(12, 5).Value = '  -0.21%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.031'
$ws.Cells.Item(13, 5).Value = '  +0.82%  '

$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.6807'
$ws.Cells.Item(14, 5).Value = '  +0.37%  '

$ws.Cells.Item(15, 2).Value = 'ShibaInu'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.00001072'
$ws.Cells.Item(15, 5).Value = '  +3.54%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '83.41'
$ws.Cells.Item(16, 5).Value = '  +0.67%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '6.174'
$ws.Cells.Item(17, 5).Value = '  +0.66%  '

$ws.Cells.Item(18, 4).Value = '29.448.81'
$ws.Cells.Item(18, 5).Value = '  +0.13%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '228.29'
$ws.Cells.Item(19, 5).Value = '  +0.21%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.40'
$ws.Cells.Item(20, 5).Value = '  -0.07%  '

$ws.Cells.Item(21, 5).Value = '  +0.06%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.418'
$ws.Cells.Item(22, 5).Value = '  -0.24%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(23, 5).Value = '  +0.05%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '157.94'
$ws.Cells.Item(24, 5).Value = '  +0.64%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.1376'
$ws.Cells.Item(25, 5).Value = '  -1.04%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.407'
$ws.Cells.Item(26, 5).Value = '  +0.72%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '17.69'
$ws.Cells.Item(27, 5).Value = '  +0.59%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.348'
$ws.Cells.Item(28, 5).Value = '  +5.69%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.464'
$ws.Cells.Item(29, 5).Value = '  +0.45%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.05670'
$ws.Cells.Item(30, 5).Value = '  +0.89%  '

$ws.Cells.Item(31, 5).Value = '  +0.53%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.026'
$ws.Cells.Item(32, 5).Value = '  +0.23%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.843'
$ws.Cells.Item(33, 5).Value = '  +0.71%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.163'
$ws.Cells.Item(34, 5).Value = '  +0.62%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7026'
$ws.Cells.Item(35, 5).Value = '  -1.17%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.581'
$ws.Cells.Item(36, 5).Value = '  -0.25%  '

$ws.Cells.Item(37, 4).Value = '1.226.80'
$ws.Cells.Item(37, 5).Value = '  -1.08%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.764'
$ws.Cells.Item(38, 5).Value = '  -0.04%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01789'
$ws.Cells.Item(39, 5).Value = '  -0.91%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.542'
$ws.Cells.Item(40, 5).Value = '  +3.59%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.9029'
$ws.Cells.Item(41, 5).Value = '  +0.39%  '

$ws.Cells.Item(43, 4).Value = '1.999.46'
$ws.Cells.Item(43, 5).Value = '  -1.37%  '

$ws.Cells.Item(44, 5).Value = '  -0.16%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '66.05'
$ws.Cells.Item(45, 5).Value = '  +0.82%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.00000000122'
$ws.Cells.Item(46, 5).Value = '  +2.53%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.170'
$ws.Cells.Item(47, 5).Value = '  +1.64%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4017'
$ws.Cells.Item(48, 5).Value = '  +0.47%  '

$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.1154'
$ws.Cells.Item(49, 5).Value = '  +3.26%  '

$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '8.976'
$ws.Cells.Item(50, 5).Value = '  +0.76%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.674'
$ws.Cells.Item(51, 5).Value = '  +0.58%  '
